# Updates the "Price" (column D) and "Volume(1h)" (column E) figures on
# Sheet1 with refreshed crypto-ranking quotes (GitHub Actions symbol-list
# refresh for Sun Jan 22 22:55:12 UTC 2023).
#
# The sheet stores these figures as plain text (e.g. "300.05", "-1.16%"),
# not as numbers/percentages, so each new value is written with a leading
# apostrophe (forces text entry, same as a user typing '301.33 into the
# cell) and then ClearFormats() strips the "quote prefix" cell style Excel
# would otherwise stamp on the cell, keeping the cell formatted exactly as
# it was before (no style index) while the stored value is still text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "'301.33"
    "E2" = "'-0.70%"
    "D3" = "'36.65"
    "E3" = "'3.02%"
    "D4" = "'4.988"
    "E4" = "'-1.86%"
    "D5" = "'0.07680"
    "E5" = "'-1.28%"
    "D6" = "'2.059"
    "E6" = "'-8.68%"
    "D7" = "'7.901"
    "E7" = "'-2.02%"
    "D8" = "'0.9188"
    "E8" = "'-1.05%"
    "D9" = "'0.09666"
    "E9" = "'3.72%"
    "E10" = "'1.33%"
    "D11" = "'0.08514"
    "E11" = "'-0.71%"
    "D12" = "'0.03524"
    "E12" = "'-5.27%"
    "D13" = "'0.09940"
    "E13" = "'0.11%"
    "D14" = "'0.001478"
    "E14" = "'0.10%"
    "D15" = "'0.005626"
    "E15" = "'-2.17%"
    "D16" = "'3.460"
    "E16" = "'-0.72%"
    "D17" = "'4.023"
    "E17" = "'-0.75%"
    "D18" = "'2.419"
    "E18" = "'10.73%"
    "D19" = "'0.3384"
    "E19" = "'-2.29%"
    "D20" = "'0.1336"
    "E20" = "'1.03%"
    "D21" = "'4.789"
    "E21" = "'5.24%"
    "D22" = "'0.2197"
    "E22" = "'-1.77%"
    "D23" = "'0.04575"
    "E23" = "'-2.07%"
    "D24" = "'0.005080"
    "E24" = "'12.12%"
    "E25" = "'-0.26%"
    "D26" = "'0.0001398"
    "E26" = "'7.30%"
    "D39" = "'0.01752"
    "E39" = "'-1.13%"
    "D40" = "'0.04622"
    "E40" = "'-2.27%"
    "D41" = "'0.007454"
    "E41" = "'-5.92%"
    "D42" = "'0.1387"
    "E42" = "'-2.00%"
    "D43" = "'0.007712"
    "E43" = "'-2.44%"
    "D44" = "'0.002247"
    "E44" = "'0.97%"
    "D45" = "'0.01034"
    "E45" = "'7.40%"
    "D46" = "'0.00006269"
    "E46" = "'1.05%"
    "D47" = "'0.00000000749"
    "E47" = "'-0.35%"
    "D48" = "'0.0005792"
    "E48" = "'-0.14%"
    "D49" = "'35.37"
    "E49" = "'511.22%"
    "D50" = "'0.001997"
    "E50" = "'-25.91%"
    "D51" = "'0.00002097"
    "E51" = "'-0.35%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.Value = $updates[$cellRef]
    $cell.ClearFormats()
}
